# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Lane Late" Naranja (Femacal de La Calera)
# right above the existing row 502, pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 502 (shifts 502.. down to 504..)
$ws.Rows.Item(502).Insert()
$ws.Rows.Item(502).Insert()

# New row 502: Lane Late / Primera
$ws.Cells.Item(502, 1).Value = 3
$ws.Cells.Item(502, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(502, 3).Value = "Coquimbo"
$ws.Cells.Item(502, 4).Value = 44509
$ws.Cells.Item(502, 5).Value = 5
$ws.Cells.Item(502, 6).Value = "Fruta"
$ws.Cells.Item(502, 7).Value = 100102
$ws.Cells.Item(502, 8).Value = "Cítricos"
$ws.Cells.Item(502, 9).Value = 100102005
$ws.Cells.Item(502, 10).Value = "Naranja"
$ws.Cells.Item(502, 11).Value = "Lane Late"
$ws.Cells.Item(502, 12).Value = "Primera"
$ws.Cells.Item(502, 13).Value = 225
$ws.Cells.Item(502, 14).Value = 4500
$ws.Cells.Item(502, 15).Value = 5000
$ws.Cells.Item(502, 16).Value = 4778
$ws.Cells.Item(502, 17).Value = "$/malla 13 kilos"
$ws.Cells.Item(502, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(502, 19).Value = 368
$ws.Cells.Item(502, 20).Value = 13

# New row 503: Lane Late / Segunda
$ws.Cells.Item(503, 1).Value = 3
$ws.Cells.Item(503, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(503, 3).Value = "Coquimbo"
$ws.Cells.Item(503, 4).Value = 44509
$ws.Cells.Item(503, 5).Value = 5
$ws.Cells.Item(503, 6).Value = "Fruta"
$ws.Cells.Item(503, 7).Value = 100102
$ws.Cells.Item(503, 8).Value = "Cítricos"
$ws.Cells.Item(503, 9).Value = 100102005
$ws.Cells.Item(503, 10).Value = "Naranja"
$ws.Cells.Item(503, 11).Value = "Lane Late"
$ws.Cells.Item(503, 12).Value = "Segunda"
$ws.Cells.Item(503, 13).Value = 220
$ws.Cells.Item(503, 14).Value = 3500
$ws.Cells.Item(503, 15).Value = 4000
$ws.Cells.Item(503, 16).Value = 3773
$ws.Cells.Item(503, 17).Value = "$/malla 13 kilos"
$ws.Cells.Item(503, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(503, 19).Value = 290
$ws.Cells.Item(503, 20).Value = 13
